# Update "想去人数" (want-to-go count) values on sheets 展览, 演出, 全部类型
# to reflect newly generated output (commit: "Update gh-pages to output
# generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 765
$ws1.Range("F4").Value = 54
$ws1.Range("F6").Value = 116
$ws1.Range("F8").Value = 116
$ws1.Range("F9").Value = 318
$ws1.Range("F10").Value = 430
$ws1.Range("F13").Value = 11345
$ws1.Range("F14").Value = 5365

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 94
$ws2.Range("F4").Value = 6

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 765
$ws4.Range("F4").Value = 54
$ws4.Range("F5").Value = 94
$ws4.Range("F8").Value = 116
$ws4.Range("F10").Value = 116
$ws4.Range("F11").Value = 318
$ws4.Range("F12").Value = 430
$ws4.Range("F15").Value = 11345
$ws4.Range("F16").Value = 6
$ws4.Range("F17").Value = 5365
